$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'BTC'
$ws.Range("C2").Value = 'Bitcoin'
$ws.Range("D2").Value = 25810
$ws.Range("E2").Value = 502572400738
$ws.Range("F2").Value = 4541272403
$ws.Range("G2").Value = -0.16589

$ws.Range("B3").Value = 'ETH'
$ws.Range("C3").Value = 'Ethereum'
$ws.Range("D3").Value = 1625.83
$ws.Range("E3").Value = 195386934022
$ws.Range("F3").Value = 3420405220
$ws.Range("G3").Value = -0.5466800000000001

$ws.Range("B4").Value = 'USDT'
$ws.Range("C4").Value = 'Tether'
$ws.Range("D4").Value = 0.999846
$ws.Range("E4").Value = 83002670998
$ws.Range("F4").Value = 10100280055
$ws.Range("G4").Value = 0.00924

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'BNB'
$ws.Range("D5").Value = 213.03
$ws.Range("E5").Value = 32749186633
$ws.Range("F5").Value = 234969881
$ws.Range("G5").Value = -0.96637

$ws.Range("B6").Value = 'XRP'
$ws.Range("C6").Value = 'XRP'
$ws.Range("D6").Value = 0.499586
$ws.Range("E6").Value = 26516032516
$ws.Range("F6").Value = 343609476
$ws.Range("G6").Value = -0.69555

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'USD Coin'
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 26154764800
$ws.Range("F7").Value = 1518882011
$ws.Range("G7").Value = 0.05348

$ws.Range("B8").Value = 'STETH'
$ws.Range("C8").Value = 'Lido Staked Ether'
$ws.Range("D8").Value = 1625.38
$ws.Range("E8").Value = 13977019653
$ws.Range("F8").Value = 4757010
$ws.Range("G8").Value = -0.54683

$ws.Range("B9").Value = 'ADA'
$ws.Range("C9").Value = 'Cardano'
$ws.Range("D9").Value = 0.250578
$ws.Range("E9").Value = 8771955317
$ws.Range("F9").Value = 95055376
$ws.Range("G9").Value = -1.63708

$ws.Range("B10").Value = 'DOGE'
$ws.Range("C10").Value = 'Dogecoin'
$ws.Range("D10").Value = 0.061855
$ws.Range("E10").Value = 8711510172
$ws.Range("F10").Value = 210802392
$ws.Range("G10").Value = -2.79266

$ws.Range("B11").Value = 'SOL'
$ws.Range("C11").Value = 'Solana'
$ws.Range("D11").Value = 18.39
$ws.Range("E11").Value = 7548025975
$ws.Range("F11").Value = 363090373
$ws.Range("G11").Value = -5.91508

$ws.Range("B12").Value = 'TRX'
$ws.Range("C12").Value = 'TRON'
$ws.Range("D12").Value = 0.078515
$ws.Range("E12").Value = 7005447037
$ws.Range("F12").Value = 127459856
$ws.Range("G12").Value = -0.50712

$ws.Range("B13").Value = 'TON'
$ws.Range("C13").Value = 'Toncoin'
$ws.Range("D13").Value = 1.75
$ws.Range("E13").Value = 6016950131
$ws.Range("F13").Value = 13201608
$ws.Range("G13").Value = -2.45465

$ws.Range("B14").Value = 'DOT'
$ws.Range("C14").Value = 'Polkadot'
$ws.Range("D14").Value = 4.16
$ws.Range("E14").Value = 5287464634
$ws.Range("F14").Value = 73693347
$ws.Range("G14").Value = -2.59101

$ws.Range("B15").Value = 'MATIC'
$ws.Range("C15").Value = 'Polygon'
$ws.Range("D15").Value = 0.5289970000000001
$ws.Range("E15").Value = 4927451356
$ws.Range("F15").Value = 128831131
$ws.Range("G15").Value = -2.4883

$ws.Range("B16").Value = 'LTC'
$ws.Range("C16").Value = 'Litecoin'
$ws.Range("D16").Value = 61.38
$ws.Range("E16").Value = 4514830144
$ws.Range("F16").Value = 260349702
$ws.Range("G16").Value = -2.18587

$ws.Range("B17").Value = 'SHIB'
$ws.Range("C17").Value = 'Shiba Inu'
$ws.Range("D17").Value = 0.00000741
$ws.Range("E17").Value = 4363166347
$ws.Range("F17").Value = 65169882
$ws.Range("G17").Value = -2.71464

$ws.Range("B18").Value = 'WBTC'
$ws.Range("C18").Value = 'Wrapped Bitcoin'
$ws.Range("D18").Value = 25827
$ws.Range("E18").Value = 4203782118
$ws.Range("F18").Value = 17790539
$ws.Range("G18").Value = -0.27365

$ws.Range("B19").Value = 'DAI'
$ws.Range("C19").Value = 'Dai'
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 3834891838
$ws.Range("F19").Value = 40445432
$ws.Range("G19").Value = 0.07126

$ws.Range("B20").Value = 'BCH'
$ws.Range("C20").Value = 'Bitcoin Cash'
$ws.Range("D20").Value = 189.27
$ws.Range("E20").Value = 3684649922
$ws.Range("F20").Value = 114084241
$ws.Range("G20").Value = -1.2484

$ws.Range("B21").Value = 'XLM'
$ws.Range("C21").Value = 'Stellar'
$ws.Range("D21").Value = 0.13151
$ws.Range("E21").Value = 3600439224
$ws.Range("F21").Value = 116742276
$ws.Range("G21").Value = -0.24153

$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'LEO Token'
$ws.Range("D22").Value = 3.85
$ws.Range("E22").Value = 3582222481
$ws.Range("F22").Value = 133997
$ws.Range("G22").Value = 2.00064

$ws.Range("B23").Value = 'AVAX'
$ws.Range("C23").Value = 'Avalanche'
$ws.Range("D23").Value = 9.56
$ws.Range("E23").Value = 3367882831
$ws.Range("F23").Value = 100070315
$ws.Range("G23").Value = -3.34149

$ws.Range("B24").Value = 'LINK'
$ws.Range("C24").Value = 'Chainlink'
$ws.Range("D24").Value = 6.09
$ws.Range("E24").Value = 3277434808
$ws.Range("F24").Value = 124532301
$ws.Range("G24").Value = -1.9651

$ws.Range("B25").Value = 'TUSD'
$ws.Range("C25").Value = 'TrueUSD'
$ws.Range("D25").Value = 0.998596
$ws.Range("E25").Value = 3182576832
$ws.Range("F25").Value = 129415266
$ws.Range("G25").Value = 0.03605

$ws.Range("B26").Value = 'UNI'
$ws.Range("C26").Value = 'Uniswap'
$ws.Range("D26").Value = 4.22
$ws.Range("E26").Value = 3178704633
$ws.Range("F26").Value = 71035227
$ws.Range("G26").Value = -2.52796

$ws.Range("B27").Value = 'BUSD'
$ws.Range("C27").Value = 'Binance USD'
$ws.Range("D27").Value = 1.001
$ws.Range("E27").Value = 2590861430
$ws.Range("F27").Value = 1171585271
$ws.Range("G27").Value = 0.07191

$ws.Range("B28").Value = 'XMR'
$ws.Range("C28").Value = 'Monero'
$ws.Range("D28").Value = 142.48
$ws.Range("E28").Value = 2584957418
$ws.Range("F28").Value = 32696345
$ws.Range("G28").Value = -0.58231

$ws.Range("B29").Value = 'OKB'
$ws.Range("C29").Value = 'OKB'
$ws.Range("D29").Value = 41.85
$ws.Range("E29").Value = 2508596267
$ws.Range("F29").Value = 1974297
$ws.Range("G29").Value = -0.82369

$ws.Range("B30").Value = 'ETC'
$ws.Range("C30").Value = 'Ethereum Classic'
$ws.Range("D30").Value = 15.12
$ws.Range("E30").Value = 2160942273
$ws.Range("F30").Value = 45265153
$ws.Range("G30").Value = -1.8132

$ws.Range("B31").Value = 'ATOM'
$ws.Range("C31").Value = 'Cosmos Hub'
$ws.Range("D31").Value = 6.76
$ws.Range("E31").Value = 1974384776
$ws.Range("F31").Value = 85204957
$ws.Range("G31").Value = -1.91456

$ws.Range("B32").Value = 'HBAR'
$ws.Range("C32").Value = 'Hedera'
$ws.Range("D32").Value = 0.04797207
$ws.Range("E32").Value = 1592166311
$ws.Range("F32").Value = 20368711
$ws.Range("G32").Value = -3.24234

$ws.Range("B33").Value = 'QNT'
$ws.Range("C33").Value = 'Quant'
$ws.Range("D33").Value = 98.04000000000001
$ws.Range("E33").Value = 1424904498
$ws.Range("F33").Value = 11371847
$ws.Range("G33").Value = -1.10383

$ws.Range("B34").Value = 'FIL'
$ws.Range("C34").Value = 'Filecoin'
$ws.Range("D34").Value = 3.14
$ws.Range("E34").Value = 1402654454
$ws.Range("F34").Value = 70471845
$ws.Range("G34").Value = -2.92886

$ws.Range("B35").Value = 'ICP'
$ws.Range("C35").Value = 'Internet Computer'
$ws.Range("D35").Value = 3.13
$ws.Range("E35").Value = 1387615616
$ws.Range("F35").Value = 16985546
$ws.Range("G35").Value = -4.33338

$ws.Range("B36").Value = 'CRO'
$ws.Range("C36").Value = 'Cronos'
$ws.Range("D36").Value = 0.052346
$ws.Range("E36").Value = 1373858278
$ws.Range("F36").Value = 5901636
$ws.Range("G36").Value = -0.76576

$ws.Range("B37").Value = 'MNT'
$ws.Range("C37").Value = 'Mantle'
$ws.Range("D37").Value = 0.411831
$ws.Range("E37").Value = 1332171118
$ws.Range("F37").Value = 8308340
$ws.Range("G37").Value = -0.46895

$ws.Range("B38").Value = 'LDO'
$ws.Range("C38").Value = 'Lido DAO'
$ws.Range("D38").Value = 1.49
$ws.Range("E38").Value = 1320996090
$ws.Range("F38").Value = 20344793
$ws.Range("G38").Value = -2.49411

$ws.Range("B39").Value = 'APT'
$ws.Range("C39").Value = 'Aptos'
$ws.Range("D39").Value = 5.26
$ws.Range("E39").Value = 1207945936
$ws.Range("F39").Value = 47304259
$ws.Range("G39").Value = -4.05173

$ws.Range("B40").Value = 'VET'
$ws.Range("C40").Value = 'VeChain'
$ws.Range("D40").Value = 0.01539834
$ws.Range("E40").Value = 1119107073
$ws.Range("F40").Value = 23096777
$ws.Range("G40").Value = -1.19616

$ws.Range("B41").Value = 'ARB'
$ws.Range("C41").Value = 'Arbitrum'
$ws.Range("D41").Value = 0.858495
$ws.Range("E41").Value = 1094266274
$ws.Range("F41").Value = 87288365
$ws.Range("G41").Value = -4.66631

$ws.Range("B42").Value = 'NEAR'
$ws.Range("C42").Value = 'NEAR Protocol'
$ws.Range("D42").Value = 1.15
$ws.Range("E42").Value = 1084854956
$ws.Range("F42").Value = 52143362
$ws.Range("G42").Value = -0.5269

$ws.Range("B43").Value = 'MKR'
$ws.Range("C43").Value = 'Maker'
$ws.Range("D43").Value = 1130.13
$ws.Range("E43").Value = 1017781504
$ws.Range("F43").Value = 39836692
$ws.Range("G43").Value = 0.42214

$ws.Range("B44").Value = 'OP'
$ws.Range("C44").Value = 'Optimism'
$ws.Range("D44").Value = 1.27
$ws.Range("E44").Value = 1010841617
$ws.Range("F44").Value = 50161313
$ws.Range("G44").Value = -2.05708

$ws.Range("B45").Value = 'RETH'
$ws.Range("C45").Value = 'Rocket Pool ETH'
$ws.Range("D45").Value = 1762.65
$ws.Range("E45").Value = 925453190
$ws.Range("F45").Value = 3994106
$ws.Range("G45").Value = -0.47508

$ws.Range("B46").Value = 'KAS'
$ws.Range("C46").Value = 'Kaspa'
$ws.Range("D46").Value = 0.03944771
$ws.Range("E46").Value = 810610232
$ws.Range("F46").Value = 6394252
$ws.Range("G46").Value = 2.30071

$ws.Range("B47").Value = 'AAVE'
$ws.Range("C47").Value = 'Aave'
$ws.Range("D47").Value = 54.89
$ws.Range("E47").Value = 797125957
$ws.Range("F47").Value = 42325523
$ws.Range("G47").Value = -2.50719

$ws.Range("B48").Value = 'GRT'
$ws.Range("C48").Value = 'The Graph'
$ws.Range("D48").Value = 0.084642
$ws.Range("E48").Value = 776118330
$ws.Range("F48").Value = 16744748
$ws.Range("G48").Value = -2.7648

$ws.Range("B49").Value = 'WBT'
$ws.Range("C49").Value = 'WhiteBIT Coin'
$ws.Range("D49").Value = 5.29
$ws.Range("E49").Value = 761992873
$ws.Range("F49").Value = 7216652
$ws.Range("G49").Value = -0.09889000000000001

$ws.Range("B50").Value = 'XDC'
$ws.Range("C50").Value = 'XDC Network'
$ws.Range("D50").Value = 0.054645
$ws.Range("E50").Value = 757674777
$ws.Range("F50").Value = 5277795
$ws.Range("G50").Value = -2.84529

$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'USDD'
$ws.Range("D51").Value = 0.997434
$ws.Range("E51").Value = 723197562
$ws.Range("F51").Value = 17825348
$ws.Range("G51").Value = 0.0483
